$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 641.25
$ws.Range("I12").Value = 509.2
$ws.Range("J12").Value = 861.3333
$ws.Range("K12").Value = 509.2
$ws.Range("L12").Value = 861.3333
$ws.Range("M12").Value = -339.2
$ws.Range("N12").Value = -1201.3333
$ws.Range("H17").Value = 2155.4092
$ws.Range("J17").Value = 2205.6667
$ws.Range("L17").Value = 6617.000100000001
$ws.Range("N17").Value = -6953.000100000001
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = 15
$ws.Range("M29").Value = 266
$ws.Range("H64").Value = 4199.6
$ws.Range("J64").Value = 3001
$ws.Range("L64").Value = 3001
$ws.Range("N64").Value = -3497
$ws.Range("H67").Value = 4199.6
$ws.Range("J67").Value = 3001
$ws.Range("L67").Value = 3001
$ws.Range("N67").Value = -4717
$ws.Range("H137").Value = 1451988.1
$ws.Range("I137").Value = 1862.5
$ws.Range("J137").Value = 2418738.5
$ws.Range("K137").Value = 5587.5
$ws.Range("L137").Value = 7256215.5
$ws.Range("M137").Value = -3037.5
$ws.Range("N137").Value = -7261315.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6672.1665
$ws.Range("I32").Value = 3223.9253
$ws.Range("K32").Value = 3223.9253
$ws.Range("M32").Value = -2936.9253
$ws.Range("H61").Value = 126780.75
$ws.Range("I61").Value = 1941.1666
$ws.Range("J61").Value = 501299.5
$ws.Range("K61").Value = 1941.1666
$ws.Range("L61").Value = 501299.5
$ws.Range("M61").Value = -1729.1666
$ws.Range("N61").Value = -501723.5
$ws.Range("H74").Value = 37536.715
$ws.Range("I74").Value = 60139.824
$ws.Range("J74").Value = 2604.6365
$ws.Range("K74").Value = 60139.824
$ws.Range("L74").Value = 2604.6365
$ws.Range("M74").Value = -59265.824
$ws.Range("N74").Value = -4352.636500000001
$ws.Range("H77").Value = 37536.715
$ws.Range("I77").Value = 60139.824
$ws.Range("J77").Value = 2604.6365
$ws.Range("K77").Value = 300699.12
$ws.Range("L77").Value = 13023.1825
$ws.Range("M77").Value = -296331.12
$ws.Range("N77").Value = -21759.1825
$ws.Range("H97").Value = 646.3684
$ws.Range("I97").Value = 755.9286
$ws.Range("K97").Value = 755.9286
$ws.Range("M97").Value = -259.9286
$ws.Range("H122").Value = 3796.4285
$ws.Range("I122").Value = 3961.2222
$ws.Range("J122").Value = 3499.8
$ws.Range("K122").Value = 11883.6666
$ws.Range("L122").Value = 10499.4
$ws.Range("M122").Value = -9433.6666
$ws.Range("N122").Value = -15399.4
$ws.Range("H132").Value = 1865.1794
$ws.Range("I132").Value = 1627.9706
$ws.Range("J132").Value = 3478.2
$ws.Range("K132").Value = 4883.9118
$ws.Range("L132").Value = 10434.6
$ws.Range("M132").Value = -2353.9118
$ws.Range("N132").Value = -15494.6
$ws.Range("H136").Value = 126780.75
$ws.Range("I136").Value = 1941.1666
$ws.Range("J136").Value = 501299.5
$ws.Range("K136").Value = 5823.4998
$ws.Range("L136").Value = 1503898.5
$ws.Range("M136").Value = -3273.4998
$ws.Range("N136").Value = -1508998.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 72726.92999999999
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H86").Value = 2396.25
$ws.Range("I86").Value = 1705
$ws.Range("K86").Value = 1705
$ws.Range("M86").Value = -582
$ws.Range("H89").Value = 2396.25
$ws.Range("I89").Value = 1705
$ws.Range("K89").Value = 8525
$ws.Range("M89").Value = -2909
$ws.Range("H134").Value = 1927.44
$ws.Range("I134").Value = 1359.5
$ws.Range("K134").Value = 4078.5
$ws.Range("M134").Value = -1543.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1594.2307
$ws.Range("I16").Value = 1171.2
$ws.Range("K16").Value = 1171.2
$ws.Range("M16").Value = -884.2
$ws.Range("H22").Value = 579.8
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -1699.5
$ws.Range("H58").Value = 2205
$ws.Range("I58").Value = 2205
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2205
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2002
$ws.Range("N58").ClearContents()
$ws.Range("H105").Value = 2948.3333
$ws.Range("I105").Value = 636
$ws.Range("J105").Value = 4600
$ws.Range("K105").Value = 636
$ws.Range("L105").Value = 4600
$ws.Range("M105").Value = 1111
$ws.Range("N105").Value = -8094
$ws.Range("H113").Value = 1594.2307
$ws.Range("I113").Value = 1171.2
$ws.Range("K113").Value = 1171.2
$ws.Range("M113").Value = 998.8
$ws.Range("H122").Value = 2928.4375
$ws.Range("I122").Value = 1959.5454
$ws.Range("K122").Value = 5878.6362
$ws.Range("M122").Value = -3428.6362
$ws.Range("H132").Value = 2781.4
$ws.Range("I132").Value = 2781.4
$ws.Range("K132").Value = 8344.200000000001
$ws.Range("M132").Value = -5814.200000000001
$ws.Range("H136").Value = 2205
$ws.Range("I136").Value = 2205
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6615
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4065
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 127.47619
$ws.Range("I2").Value = 104.30769
$ws.Range("J2").Value = 165.125
$ws.Range("K2").Value = 625.84614
$ws.Range("L2").Value = 990.75
$ws.Range("M2").Value = -512.84614
$ws.Range("N2").Value = -1216.75
$ws.Range("H4").Value = 114253860
$ws.Range("I4").Value = 133181160
$ws.Range("J4").Value = 690000
$ws.Range("K4").Value = 399543480
$ws.Range("L4").Value = 2070000
$ws.Range("M4").Value = -399543368
$ws.Range("N4").Value = -2070224
$ws.Range("H6").Value = 183152.27
$ws.Range("I6").Value = 223841.67
$ws.Range("K6").Value = 671525.01
$ws.Range("M6").Value = -671412.01
$ws.Range("H35").Value = 10486.667
$ws.Range("I35").Value = 15230
$ws.Range("K35").Value = 45690
$ws.Range("M35").Value = -45402

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 28149.834
$ws.Range("I58").Value = 26599.666
$ws.Range("J58").Value = 29700
$ws.Range("K58").Value = 26599.666
$ws.Range("L58").Value = 29700
$ws.Range("M58").Value = -26322.666
$ws.Range("N58").Value = -30254
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H113").Value = 3705788.2
$ws.Range("I113").Value = 1166.6666
$ws.Range("J113").Value = 5558099
$ws.Range("K113").Value = 1166.6666
$ws.Range("L113").Value = 5558099
$ws.Range("M113").Value = 1003.3334
$ws.Range("N113").Value = -5562439
$ws.Range("H132").Value = 3133.2322
$ws.Range("I132").Value = 2444.8635
$ws.Range("K132").Value = 7334.5905
$ws.Range("M132").Value = -4804.5905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I16").Value = 1490.8889
$ws.Range("J16").Value = 1149.5
$ws.Range("K16").Value = 1490.8889
$ws.Range("L16").Value = 1149.5
$ws.Range("M16").Value = -1320.8889
$ws.Range("N16").Value = -1489.5
$ws.Range("H22").Value = 5840.6665
$ws.Range("I22").Value = 894.875
$ws.Range("J22").Value = 8884.23
$ws.Range("K22").Value = 894.875
$ws.Range("L22").Value = 8884.23
$ws.Range("M22").Value = -599.875
$ws.Range("N22").Value = -9474.23
$ws.Range("H27").Value = 5840.6665
$ws.Range("I27").Value = 894.875
$ws.Range("J27").Value = 8884.23
$ws.Range("K27").Value = 894.875
$ws.Range("L27").Value = 8884.23
$ws.Range("M27").Value = -787.875
$ws.Range("N27").Value = -9098.23
$ws.Range("H61").Value = 1880.8889
$ws.Range("I61").Value = 1877.875
$ws.Range("K61").Value = 1877.875
$ws.Range("M61").Value = -1675.875
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 60000
$ws.Range("M65").Value = -56880
$ws.Range("H113").Value = 1880.8889
$ws.Range("I113").Value = 1877.875
$ws.Range("K113").Value = 1877.875
$ws.Range("M113").Value = 292.125
$ws.Range("H122").Value = 16696801
$ws.Range("I122").Value = 48214.57
$ws.Range("J122").Value = 40004820
$ws.Range("K122").Value = 144643.71
$ws.Range("L122").Value = 120014460
$ws.Range("M122").Value = -142193.71
$ws.Range("N122").Value = -120019360
$ws.Range("H132").Value = 8943.333000000001
$ws.Range("I132").Value = 10154.235
$ws.Range("K132").Value = 30462.705
$ws.Range("M132").Value = -27932.705
$ws.Range("H136").Value = 6828.7827
$ws.Range("J136").Value = 6453.6924
$ws.Range("L136").Value = 19361.0772
$ws.Range("N136").Value = -24461.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 891.4286
$ws.Range("I113").Value = 806.2308
$ws.Range("K113").Value = 2418.6924
$ws.Range("M113").Value = -248.6923999999999
$ws.Range("I122").Value = 2611.182
$ws.Range("J122").Value = 3293.4
$ws.Range("K122").Value = 7833.545999999999
$ws.Range("L122").Value = 9880.200000000001
$ws.Range("M122").Value = -5383.545999999999
$ws.Range("N122").Value = -14780.2
$ws.Range("H132").Value = 1674101.1
$ws.Range("I132").Value = 1660.4286
$ws.Range("K132").Value = 4981.2858
$ws.Range("M132").Value = -2451.2858
$ws.Range("H136").Value = 1497.25
$ws.Range("I136").Value = 1452.5555
$ws.Range("J136").Value = 1899.5
$ws.Range("K136").Value = 4357.666499999999
$ws.Range("L136").Value = 5698.5
$ws.Range("M136").Value = -1807.666499999999
$ws.Range("N136").Value = -10798.5
